# Refresh the crypto price/volume columns (D, E) with the latest snapshot.
# Values that look like plain numbers (e.g. "0.9620") are prefixed with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inline-string cells) instead of silently parsing them into
# numbers and dropping significant trailing/leading zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.551.20"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "1.473.16"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'0.9620"
$ws.Range("E5").Value = "  +5.09%  "
$ws.Range("D6").Value = "'277.11"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").Value = "'0.3594"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").Value = "'0.3087"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "'1.084"
$ws.Range("E9").Value = "  +5.91%  "
$ws.Range("D10").Value = "'39.51"
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").Value = "'18.18"
$ws.Range("E14").Value = "  +4.06%  "
$ws.Range("D15").Value = "'6.170"
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "'0.9613"
$ws.Range("E16").Value = "  +2.75%  "
$ws.Range("D17").Value = "'0.00001023"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "1.470.96"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "'0.05945"
$ws.Range("D20").Value = "'69.00"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").Value = "'5.495"
$ws.Range("E21").Value = "  +2.12%  "
$ws.Range("D22").Value = "'14.59"
$ws.Range("E22").Value = "  +2.33%  "
$ws.Range("D23").Value = "'11.27"
$ws.Range("E23").Value = "  +4.02%  "
$ws.Range("D24").Value = "'2.274"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "20.546.75"
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "'144.21"
$ws.Range("E26").Value = "  +4.45%  "
$ws.Range("D27").Value = "'2.119"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "'17.18"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").Value = "1.632.75"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").Value = "'113.81"
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("D31").Value = "'3.891"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'0.8099"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").Value = "'4.942"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "'0.08002"
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("D35").Value = "'1.227"
$ws.Range("E35").Value = "  +8.91%  "
$ws.Range("D36").Value = "'1.467"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "'0.05799"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").Value = "'4.728"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").Value = "'0.02050"
$ws.Range("E39").Value = "  +3.30%  "
$ws.Range("D40").Value = "'0.9622"
$ws.Range("E40").Value = "  +3.72%  "
$ws.Range("D41").Value = "'10.40"
$ws.Range("E41").Value = "  +2.19%  "
$ws.Range("D42").Value = "'0.1874"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").Value = "'7.404"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("D44").Value = "'0.5274"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").Value = "'3.518"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "'12.20"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").Value = "'119.25"
$ws.Range("E47").Value = "  +2.03%  "
$ws.Range("D48").Value = "'0.5204"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").Value = "'1.816"
$ws.Range("E49").Value = "  +4.94%  "
$ws.Range("D50").Value = "'0.06453"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").Value = "'0.9936"
$ws.Range("E51").Value = "  +0.61%  "
